$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 41665.5135
$ws.Range("A19").Value = 39727.5055
$ws.Range("A20").Value = 39165.7675
$ws.Range("A21").Value = 37817.495
$ws.Range("A22").Value = 38010.227
$ws.Range("A23").Value = 38285.2065
